# Updated cryptos list (price/volume refresh + two coin-row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.948.58'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.291.85'
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.14'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.641'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.61'
$ws.Range("E7").Value = '  +2.04%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.13'
$ws.Range("E10").Value = '  -4.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0978'
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.15'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.44'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").Value = '2.633.21'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.38'
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.871'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = '2.296.07'
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").Value = '42.846.04'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.30'
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.66'
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.25'
$ws.Range("E23").Value = '  +6.68%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.58'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.90'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.60'
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.42'
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.64'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.19'
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.94'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.08'
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.38'
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0833'
$ws.Range("E34").Value = '  +5.75%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.12'
$ws.Range("E36").Value = '  +7.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.127'
$ws.Range("E37").Value = '  +1.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.60'
$ws.Range("E38").Value = '  +11.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.80'
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0309'
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.88'
$ws.Range("E41").Value = '  +11.48%  '
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.91'
$ws.Range("E43").Value = '  -1.82%  '
$ws.Range("E44").Value = '  +8.00%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.14'
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.03'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.27'
$ws.Range("E49").Value = '  +9.25%  '
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("E51").Value = '  -1.80%  '
